$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the five new task/header labels across row 2 (columns C-G).
# Typing these values creates new shared-string entries (indices 7-11),
# matching the sharedStrings.xml additions in the diff.
$ws.Range("C2").Value = "Client/appli android"
$ws.Range("D2").Value = "Serveur/raspberry pi"
$ws.Range("E2").Value = "Serveur : Envoie de la vidéo"
$ws.Range("F2").Value = "Client : réception/affichage de la vidéo"
$ws.Range("G2").Value = "Client : commande des roues"

# Widen the new columns so the headers are readable (mirrors the <cols>
# block added to the sheet). The COM layer quantizes ColumnWidth to whole
# pixels, so these values are chosen to land on the same rounded width
# Excel itself produced.
$ws.Columns.Item(3).ColumnWidth = 23.8325
$ws.Columns.Item(4).ColumnWidth = 20.3325
$ws.Columns.Item(5).ColumnWidth = 26.1675
$ws.Columns.Item(6).ColumnWidth = 36.0
$ws.Columns.Item(7).ColumnWidth = 26.6675

# Move the active selection off the header area, to H2, as in the diff.
$ws.Range("H2").Select()
